$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "sdfgfsdg"
$ws.Range("C8").Value = "dfg"
$ws.Range("C9").Value = "dfg"
$ws.Range("G10").Value = "df"
$ws.Range("G7").Value = "er"
$ws.Range("D7").Value = "t45645645"

[void]$ws.Range("F9").Select()
